$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ===================================================================
# Update coin identity (Name / Link) for rows whose ranking order
# changed position in the source feed.
# ===================================================================
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("B46").Value = "BitcoinSV"
$ws.Range("C46").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"

# ===================================================================
# Update Price (column D). Values are stored as plain text in the
# workbook (e.g. "320.55"), so where the text would otherwise be
# auto-recognised by Excel as a number we prefix it with a leading
# apostrophe to force text entry, matching the original formatting.
# ===================================================================
$ws.Range("D2").Value = "44.473.07"
$ws.Range("D3").Value = "2.275.21"
$ws.Range("D5").Value = "'320.55"
$ws.Range("D6").Value = "'105.82"
$ws.Range("D9").Value = "'0.571"
$ws.Range("D10").Value = "'38.66"
$ws.Range("D11").Value = "'0.0843"
$ws.Range("D12").Value = "'7.88"
$ws.Range("D14").Value = "2.624.07"
$ws.Range("D15").Value = "'0.883"
$ws.Range("D16").Value = "'14.58"
$ws.Range("D17").Value = "2.282.95"
$ws.Range("D18").Value = "44.289.74"
$ws.Range("D19").Value = "'14.02"
$ws.Range("D20").Value = "'0.0000100"
$ws.Range("D21").Value = "'6.53"
$ws.Range("D22").Value = "'66.38"
$ws.Range("D24").Value = "'238.67"
$ws.Range("D25").Value = "'2.20"
$ws.Range("D27").Value = "'10.17"
$ws.Range("D29").Value = "'38.24"
$ws.Range("D30").Value = "'6.51"
$ws.Range("D31").Value = "'163.83"
$ws.Range("D32").Value = "'20.61"
$ws.Range("D36").Value = "'0.116"
$ws.Range("D37").Value = "'3.18"
$ws.Range("D39").Value = "'3.97"
$ws.Range("D40").Value = "'4.45"
$ws.Range("D41").Value = "'15.54"
$ws.Range("D42").Value = "'0.0328"
$ws.Range("D44").Value = "1.779.16"
$ws.Range("D45").Value = "'0.208"
$ws.Range("D46").Value = "'86.71"
$ws.Range("D47").Value = "'5.49"
$ws.Range("D48").Value = "'60.23"
$ws.Range("D49").Value = "'74.89"
$ws.Range("D50").Value = "'104.35"
$ws.Range("D51").Value = "'8.66"

# ===================================================================
# Update Volume(1h) (column E). Values keep the two-space padding
# on both sides exactly like the source data.
# ===================================================================
$ws.Range("E2").Value = "  +3.51%  "
$ws.Range("E3").Value = "  +2.17%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("E6").Value = "  +5.55%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  +1.34%  "
$ws.Range("E10").Value = "  +3.34%  "
$ws.Range("E11").Value = "  +1.38%  "
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("E14").Value = "  +2.36%  "
$ws.Range("E15").Value = "  +2.20%  "
$ws.Range("E16").Value = "  +1.99%  "
$ws.Range("E17").Value = "  +2.67%  "
$ws.Range("E18").Value = "  +3.35%  "
$ws.Range("E19").Value = "  -6.11%  "
$ws.Range("E20").Value = "  +3.80%  "
$ws.Range("E21").Value = "  +1.17%  "
$ws.Range("E22").Value = "  +1.07%  "
$ws.Range("E23").Value = "  +1.26%  "
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("E25").Value = "  +2.34%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("E27").Value = "  +1.07%  "
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("E29").Value = "  +12.03%  "
$ws.Range("E30").Value = "  +1.85%  "
$ws.Range("E31").Value = "  +4.60%  "
$ws.Range("E32").Value = "  +0.42%  "
$ws.Range("E33").Value = "  -2.31%  "
$ws.Range("E34").Value = "  -0.87%  "
$ws.Range("E35").Value = "  +3.33%  "
$ws.Range("E36").Value = "  +10.40%  "
$ws.Range("E37").Value = "  -1.11%  "
$ws.Range("E38").Value = "  -0.64%  "
$ws.Range("E39").Value = "  +1.23%  "
$ws.Range("E40").Value = "  -0.43%  "
$ws.Range("E41").Value = "  +23.57%  "
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("E43").Value = "  +0.25%  "
$ws.Range("E44").Value = "  -8.13%  "
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("E46").Value = "  -2.52%  "
$ws.Range("E47").Value = "  +1.86%  "
$ws.Range("E48").Value = "  -0.50%  "
$ws.Range("E49").Value = "  -2.62%  "
$ws.Range("E50").Value = "  +1.02%  "
$ws.Range("E51").Value = "  +0.88%  "

# ===================================================================
# Reset style on Price cells that needed the text-forcing apostrophe
# so the quote-prefix formatting flag is not left behind.
# ===================================================================
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
